$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.757.67'
$ws.Range("E2").Value = '  +0.18%  '

# Row 3
$ws.Range("D3").Value = '2.087.95'
$ws.Range("E3").Value = '  +0.70%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.637'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.35%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.16'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.20%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '

# Row 10
$ws.Range("E10").Value = '  -0.02%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.108'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.12%  '

# Row 13
$ws.Range("D13").Value = '2.395.09'
$ws.Range("E13").Value = '  +0.71%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.778'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.14%  '

# Row 16
$ws.Range("E16").Value = '  +1.34%  '

# Row 17
$ws.Range("D17").Value = '2.081.90'
$ws.Range("E17").Value = '  +0.90%  '

# Row 18
$ws.Range("D18").Value = '37.774.53'
$ws.Range("E18").Value = '  +0.37%  '

# Row 19
$ws.Range("E19").Value = '  -1.25%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.15%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0835'
$ws.Range("E21").Value = '  +0.48%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.79%  '

# Row 23
$ws.Range("E23").Value = '  -0.07%  '

# Row 24
$ws.Range("E24").Value = '  -1.11%  '

# Row 25
$ws.Range("E25").Value = '  +0.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.39%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.89%  '

# Row 28
$ws.Range("E28").Value = '  -4.04%  '

# Row 29
$ws.Range("E29").Value = '  +0.70%  '

# Row 30
$ws.Range("E30").Value = '  -0.34%  '

# Row 31
$ws.Range("E31").Value = '  +1.74%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.44%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0635'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.15%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.69%  '

# Row 35
$ws.Range("E35").Value = '  +1.90%  '

# Row 36
$ws.Range("E36").Value = '  -0.06%  '

# Row 37
$ws.Range("E37").Value = '  -1.58%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.04%  '

# Row 39
$ws.Range("E39").Value = '  +1.27%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0234'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.11%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.25'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.44%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0963'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.44%  '

# Row 43
$ws.Range("E43").Value = '  +1.12%  '

# Row 44
$ws.Range("E44").Value = '  +3.02%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.47%  '

# Row 46
$ws.Range("D46").Value = '1.451.43'
$ws.Range("E46").Value = '  -0.22%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.51%  '

# Row 49
$ws.Range("E49").Value = '  -2.11%  '

# Row 50
$ws.Range("E50").Value = '  -1.94%  '

# Row 51
$ws.Range("D51").Value = '2.278.86'
$ws.Range("E51").Value = '  +0.66%  '
